$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Notes Master date placeholder: 02/02/2024 -> 11/03/2024
# -----------------------------------------------------------------
$nm = $p.NotesMaster
$dateShape = $nm.Shapes.Item(2)
$dateShape.TextFrame.TextRange.Text = "11/03/2024"

# -----------------------------------------------------------------
# 2) Slide 17: "Unreliable connectionless data transmission. " ->
#    "Finds the correct routing path, handles congestions and
#    quality of service."
# -----------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$sh17 = $s17.Shapes.Item(3)
$tr17 = $sh17.TextFrame.TextRange
$para17 = $tr17.Paragraphs(1, 1)
$para17.Text = "Finds the correct routing path, handles congestions and quality of service."

# -----------------------------------------------------------------
# 3) Slide 18 title: "The TCP/IP reference model (revised)" ->
#    "The TCP/IP reference model"
# -----------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$sh18 = $s18.Shapes.Item(1)
$tr18 = $sh18.TextFrame.TextRange
$sel18 = $tr18.Characters(21, 16)
$sel18.Text = " model"

# -----------------------------------------------------------------
# 4) Slide 22: "The request-response protocol for fetching pages..."
#    -> "The request-response application protocol for fetching
#    pages..."
# -----------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$sh22 = $s22.Shapes.Item(2)
$tr22 = $sh22.TextFrame.TextRange
$para22 = $tr22.Paragraphs(5, 1)
$run22 = $tr22.Characters($para22.Start, 108)
$run22.Text = "The request-response application protocol for fetching pages is a simple text-based protocol that runs over TCP, called "

# -----------------------------------------------------------------
# 5) Slide 26: "There are 250 top-level domains, ..." ->
#    "There are top-level domains, ..."
# -----------------------------------------------------------------
$s26 = $p.Slides.Item(26)
$sh26 = $s26.Shapes.Item(2)
$tr26 = $sh26.TextFrame.TextRange
$para26 = $tr26.Paragraphs(3, 1)
$run26 = $tr26.Characters($para26.Start, 77)
$run26.Text = "There are top-level domains, each divided into subdomains according to a "

# -----------------------------------------------------------------
# 6) Remove click-triggered animations / timing from slide 3 and
#    slide 25 (the <p:timing> element disappears entirely once all
#    effects are removed from the slide's main animation sequence).
# -----------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$seq3 = $s3.TimeLine.MainSequence
for ($i = $seq3.Count; $i -ge 1; $i--) {
    $seq3.Item($i).Delete()
}

$s25 = $p.Slides.Item(25)
$seq25 = $s25.TimeLine.MainSequence
for ($i = $seq25.Count; $i -ge 1; $i--) {
    $seq25.Item($i).Delete()
}
